$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date values in column B (rows 2-8) to remove the time component
# (shift each date back to the start of its correct day, 12 days earlier)
$ws.Range("B2").Value = 43599
$ws.Range("B3").Value = 43600
$ws.Range("B4").Value = 43601
$ws.Range("B5").Value = 43602
$ws.Range("B6").Value = 43603
$ws.Range("B7").Value = 43604
$ws.Range("B8").Value = 43605
